$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "RESOURCE_COLLECT_WITH"
$ws.Range("B17").Value = "Collect with"
$ws.Range("C17").Value = "Récolter avec"

$ws.Range("D17").Select()
